$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates exactly as described by the commit diff.
# Column D holds numeric-looking values that must remain TEXT (the source
# file stores them as inline strings, t="inlineStr"), so they are written
# via Formula with a leading apostrophe (') to force text interpretation
# and avoid Excel silently re-typing them as floating-point numbers.

$ws.Range("D2").Formula = "'242.39"
$ws.Range("D3").Formula = "'21.70"
$ws.Range("D4").Formula = "'5.376"
$ws.Range("D5").Formula = "'0.05692"
$ws.Range("D6").Formula = "'3.416"
$ws.Range("D7").Formula = "'6.290"
$ws.Range("D8").Formula = "'0.8063"
$ws.Range("D9").Formula = "'0.8373"
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Formula = "'0.01075"
$ws.Range("E10").Value = '9OneONEBestin24h'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Formula = "'0.1424"
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Formula = "'0.07267"
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Formula = "'0.03052"
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Formula = "'0.03147"
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Formula = "'0.09346"
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Formula = "'3.920"
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Formula = "'0.001587"
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Formula = "'0.04809"
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("D19").Formula = "'0.006347"
$ws.Range("D20").Formula = "'0.004060"
$ws.Range("D21").Formula = "'0.0009953"
$ws.Range("D22").Formula = "'0.0001501"
$ws.Range("D23").Formula = "'3.716"
$ws.Range("D24").Formula = "'2.153"
$ws.Range("D25").Formula = "'0.3261"
$ws.Range("D26").Formula = "'0.1300"
$ws.Range("D27").Formula = "'0.0004001"
$ws.Range("D41").Formula = "'0.006693"
$ws.Range("D42").Formula = "'0.1049"
$ws.Range("D43").Formula = "'0.002675"
$ws.Range("D44").Formula = "'0.006586"
$ws.Range("D45").Formula = "'0.00005599"
$ws.Range("D47").Formula = "'0.5802"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("D48").Formula = "'0.1424"
$ws.Range("D50").Formula = "'0.01010"
